# Applies the "Updated cryptos list" data refresh to Sheet1.
# Prices in column D that look like plain numbers are written with a leading
# apostrophe (the standard Excel text-entry prefix) so they stay text cells
# instead of being auto-converted to numbers, matching the original formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (D2,E2)
$ws.Range('D2').Value = '70.561.91'
$ws.Range('E2').Value = '  +2.11%  '

# Row 3 (D3,E3)
$ws.Range('D3').Value = '3.562.47'
$ws.Range('E3').Value = '  +1.31%  '

# Row 4 (D4,E4)
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.04%  '

# Row 5 (D5,E5)
$ws.Range('D5').Value = '''615.17'
$ws.Range('E5').Value = '  +6.33%  '

# Row 6 (D6,E6)
$ws.Range('D6').Value = '''172.19'
$ws.Range('E6').Value = '  +0.39%  '

# Row 7 (D7,E7)
$ws.Range('D7').Value = '''0.619'
$ws.Range('E7').Value = '  +1.72%  '

# Row 8 (D8,E8)
$ws.Range('D8').Value = '3.556.49'
$ws.Range('E8').Value = '  +1.27%  '

# Row 9 (E9)
$ws.Range('E9').Value = '  -0.05%  '

# Row 10 (E10)
$ws.Range('E10').Value = '  +4.06%  '

# Row 11 (D11,E11)
$ws.Range('D11').Value = '''7.19'
$ws.Range('E11').Value = '  +9.22%  '

# Row 12 (D12,E12)
$ws.Range('D12').Value = '''0.585'
$ws.Range('E12').Value = '  +0.61%  '

# Row 13 (D13,E13)
$ws.Range('D13').Value = '''46.55'
$ws.Range('E13').Value = '  -0.85%  '

# Row 14 (D14,E14)
$ws.Range('D14').Value = '''0.0000276'
$ws.Range('E14').Value = '  +1.25%  '

# Row 15 (D15,E15)
$ws.Range('D15').Value = '4.131.37'
$ws.Range('E15').Value = '  +1.02%  '

# Row 16 (D16,E16)
$ws.Range('D16').Value = '''8.36'
$ws.Range('E16').Value = '  -2.10%  '

# Row 17 (D17,E17)
$ws.Range('D17').Value = '''615.62'
$ws.Range('E17').Value = '  -1.11%  '

# Row 18 (D18,E18)
$ws.Range('D18').Value = '3.572.96'
$ws.Range('E18').Value = '  +1.62%  '

# Row 19 (D19,E19)
$ws.Range('D19').Value = '70.591.28'
$ws.Range('E19').Value = '  +2.19%  '

# Row 20 (E20)
$ws.Range('E20').Value = '  -2.20%  '

# Row 21 (D21,E21)
$ws.Range('D21').Value = '''17.40'
$ws.Range('E21').Value = '  -0.12%  '

# Row 22 (D22,E22)
$ws.Range('D22').Value = '''0.881'
$ws.Range('E22').Value = '  -0.32%  '

# Row 23 (D23,E23)
$ws.Range('D23').Value = '''9.40'
$ws.Range('E23').Value = '  -15.84%  '

# Row 24 (D24,E24)
$ws.Range('D24').Value = '''15.74'
$ws.Range('E24').Value = '  -1.26%  '

# Row 25 (D25,E25)
$ws.Range('D25').Value = '''96.65'
$ws.Range('E25').Value = '  -0.87%  '

# Row 26 (E26)
$ws.Range('E26').Value = '  +1.09%  '

# Row 27 (E27)
$ws.Range('E27').Value = '  -0.08%  '

# Row 28 (E28)
$ws.Range('E28').Value = '  -1.16%  '

# Row 29 (D29,E29)
$ws.Range('D29').Value = '''33.54'
$ws.Range('E29').Value = '  +2.78%  '

# Row 30 (D30,E30)
$ws.Range('D30').Value = '''9.03'
$ws.Range('E30').Value = '  -3.35%  '

# Row 31 (D31,E31)
$ws.Range('D31').Value = '''8.49'
$ws.Range('E31').Value = '  -0.61%  '

# Row 32 (E32)
$ws.Range('E32').Value = '  -3.81%  '

# Row 33 (E33)
$ws.Range('E33').Value = '  -1.27%  '

# Row 34 (E34)
$ws.Range('E34').Value = '  -0.76%  '

# Row 35 (D35,E35)
$ws.Range('D35').Value = '''574.48'
$ws.Range('E35').Value = '  -9.57%  '

# Row 36 (B36,C36,D36,E36)
$ws.Range('B36').Value = 'dogwifhat'
$ws.Range('C36').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D36').Value = '''3.62'
$ws.Range('E36').Value = '  +5.78%  '

# Row 37 (B37,C37,D37,E37)
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value = '''0.101'
$ws.Range('E37').Value = '  -1.51%  '

# Row 38 (D38,E38)
$ws.Range('D38').Value = '''10.82'
$ws.Range('E38').Value = '  +0.74%  '

# Row 39 (D39,E39)
$ws.Range('D39').Value = '''57.42'
$ws.Range('E39').Value = '  +1.31%  '

# Row 40 (D40,E40)
$ws.Range('D40').Value = '''0.0469'
$ws.Range('E40').Value = '  +4.89%  '

# Row 41 (E41)
$ws.Range('E41').Value = '  +0.11%  '

# Row 42 (E42)
$ws.Range('E42').Value = '  +3.66%  '

# Row 43 (D43,E43)
$ws.Range('D43').Value = '3.380.47'
$ws.Range('E43').Value = '  +0.24%  '

# Row 44 (D44,E44)
$ws.Range('D44').Value = '''0.319'
$ws.Range('E44').Value = '  -2.49%  '

# Row 45 (D45,E45)
$ws.Range('D45').Value = '''32.96'
$ws.Range('E45').Value = '  +0.12%  '

# Row 46 (D46,E46)
$ws.Range('D46').Value = '''2.96'
$ws.Range('E46').Value = '  +7.10%  '

# Row 47 (D47)
$ws.Range('D47').Value = '0.0₃0700'

# Row 48 (E48)
$ws.Range('E48').Value = '  +1.88%  '

# Row 49 (E49)
$ws.Range('E49').Value = '  +0.19%  '

# Row 50 (D50,E50)
$ws.Range('D50').Value = '''133.55'
$ws.Range('E50').Value = '  +0.61%  '
